$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "Mon Feb 24 23:06:12 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 23:06:27 EST 2025"
$ws.Range("B4").Value = "Mon Feb 24 23:06:42 EST 2025"
